# FIX — CERScore como inverso da distância. Recalculados valores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..12 (columns A..F). Row 1 header is unchanged.
$data = @(
    @("KAOGExp",     1,    1.09,               0.08445402331247004, 0.06138967062101526, 0.07850434620751935),
    @("face-knn",    1,    10.01534891301479,  47.91684236268697,   125.2692418264162,   137.3146930923626),
    @("dice",        1,    40.18008337333334,  84.01587375865105,   223.0683833392125,   113.6043565164729),
    @("cchvae",      0.93, 7.400653358181819,  45.78300191229781,   145.4949370105641,   145.1906391260749),
    @("clue",        0.93, 7.351674986727273,  48.42747684783367,   173.387140174368,    147.6900389971227),
    @("cruds",       1,    8.78503184190909,   56.51866744707791,   188.3346629580538,   184.017295453266),
    @("cem",         1,    11.53622930775,     168.3190653623361,   888.6523580938148,   192.6099408241382),
    @("wachter",     1,    5.085376672134199,  238093391.6897064,   65733031787453670,   443806747.7886946),
    @("face-epsilon",1,    9.898980877504691,  47.69467388637548,   128.240620193113,    137.6069041711729),
    @("gs",          1,    5.532854262888889,  685.255974376774,    874515.9442734839,   2808.649212049501),
    @("ar",          0.05, 0,                  0,                   0,                   0)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $rowIndex++
}

# The old sheet had 14 data rows (2..14); the new one only has 12 (2..12).
# Delete the two now-unused trailing rows (13 and 14).
$ws.Range("A13:F14").Delete()
